$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 6787
    $ws.Range("F3").Value = 49
    $ws.Range("F4").Value = 196
    $ws.Range("F5").Value = 1060
    $ws.Range("F6").Value = 150
}
